$wb = $excel.ActiveWorkbook

# --- Sheet: Rentometer ---
$ws1 = $wb.Worksheets.Item("Rentometer")

# B17: quickview_url (display text only; underlying hyperlink target is left unchanged)
$ws1.Range("B17").Value = "https://www.rentometer.com/analysis/3-bed/317-newell-st-barberton-oh-44203/-LFNYcE-hBs/quickview"

# B18: credits_remaining
$ws1.Range("B18").Value = 1955

# B19: token
$ws1.Range("B19").Value = "-LFNYcE-hBs"

# B20: links
$ws1.Range("B20").Value = "[{'rel': 'request pro report', 'href': 'https://www.rentometer.com/api/v1/request_pro_report?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=-LFNYcE-hBs'}, {'rel': 'nearby comps', 'href': 'https://www.rentometer.com/api/v1/nearby_comps?api_key=fHSGZM7POi6V5ZPR0w4CXA&token=-LFNYcE-hBs'}]"

# --- Sheet: Zillow ---
$ws2 = $wb.Worksheets.Item("Zillow")

$ws2.Range("B3").Value = 1504
$ws2.Range("B4").Value = 876
$ws2.Range("B5").Value = 1550
$ws2.Range("B8").Value = 1365.25
$ws2.Range("B9").Value = 1458.75
$ws2.Range("B10").Value = 1420

# --- Sheet: rentometer_zillow_user_avg_est ---
$ws3 = $wb.Worksheets.Item("rentometer_zillow_user_avg_est")

$ws3.Range("B1").Value = 1389
$ws3.Range("B2").Value = 1383
$ws3.Range("B3").Value = 1284.625
$ws3.Range("B4").Value = 1372.875
